$d = $word.ActiveDocument

$replacements = @(
    @("11×40=440", "71×25=1775"),
    @("26×82=2132", "75×38=2850"),
    @("16×65=1040", "40×39=1560"),
    @("29×70=2030", "77×63=4851"),
    @("84×35=2940", "89×37=3293"),
    @("66×77=5082", "36×22=792"),
    @("11×97=1067", "74×71=5254"),
    @("43×91=3913", "52×91=4732"),
    @("52×77=4004", "48×22=1056"),
    @("81×21=1701", "59×73=4307"),
    @("63×16=1008", "44×82=3608"),
    @("83×44=3652", "57×85=4845"),
    @("67×71=4757", "69×39=2691"),
    @("44×20=880", "26×68=1768"),
    @("40×91=3640", "12×78=936"),
    @("52×82=4264", "29×89=2581"),
    @("51×15=765", "32×30=960"),
    @("86×19=1634", "90×37=3330"),
    @("61×94=5734", "36×84=3024"),
    @("85×72=6120", "63×94=5922"),
    @("38×61=2318", "19×18=342"),
    @("80×17=1360", "38×69=2622"),
    @("65×13=845", "71×98=6958"),
    @("89×25=2225", "52×76=3952"),
    @("96×69=6624", "36×85=3060")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
